$d = $word.ActiveDocument
$wdParagraph = 4

# ---------------------------------------------------------------------
# 1) Remove the "License Information" Heading2 paragraph entirely.
# ---------------------------------------------------------------------
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("License Information", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $rng1.Expand($wdParagraph) | Out-Null
    $rng1.Delete()
}

# ---------------------------------------------------------------------
# 2) Rewrite the resource-description paragraph (the one that used to
#    start with the bold "Maneno Muhimu (Biblica)" run) and delete the
#    immediately-following "This PDF version..." paragraph by folding
#    both into a single paragraph replacement.
# ---------------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("Maneno Muhimu (Biblica) (Swahili) is based on", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $rng2.Expand($wdParagraph) | Out-Null
    $pStart = $rng2.Start
    $pEnd = $rng2.End

    # Clear all paragraph content but keep the paragraph mark.
    $body = $d.Range($pStart, $pEnd - 1)
    $body.Text = ""

    # --- Run 1: bold title ------------------------------------------------
    $run1Text = "Biblica Study Notes (Key Terms)"
    $insertPoint = $d.Range($pStart, $pStart)
    $insertPoint.InsertAfter($run1Text)
    $run1End = $pStart + $run1Text.Length
    $run1Range = $d.Range($pStart, $run1End)
    $run1Range.Font.Bold = 1

    # --- Run 2: license description ----------------------------------------
    $run2Text = " © 2023 Biblica Inc. Released under CC BY-SA 4.0 license. "
    $ip2 = $d.Range($run1End, $run1End)
    $ip2.InsertAfter($run2Text)
    $run2End = $run1End + $run2Text.Length
    $run2Range = $d.Range($run1End, $run2End)
    $run2Range.Font.Bold = 0

    # --- Run 3: "Biblica Study Notes" ---------------------------------------
    $run3Text = "Biblica Study Notes"
    $ip3 = $d.Range($run2End, $run2End)
    $ip3.InsertAfter($run3Text)
    $run3End = $run2End + $run3Text.Length
    $run3Range = $d.Range($run2End, $run3End)
    $run3Range.Font.Bold = 0

    # --- Run 4: adapted-languages sentence ----------------------------------
    $run4Text = " has been adapted in the following languages: Tok Pisin, Arabic (عربي), French (Français), Hindi (हिंदी), Indonesian (Bahasa Indonesia), Portuguese (Português), Russian (Русский), Spanish (Español), Swahili (Kiswahili), and Simplified Chinese (简体中文)from Biblica Study Notes © 2023 Biblica Inc. Released under CC BY-SA 4.0 license by Mission Mutual."
    $ip4 = $d.Range($run3End, $run3End)
    $ip4.InsertAfter($run4Text)
    $run4End = $run3End + $run4Text.Length
    $run4Range = $d.Range($run3End, $run4End)
    $run4Range.Font.Bold = 0
}

# ---------------------------------------------------------------------
# 3) Remove the "This PDF version is provided under the same license."
#    paragraph entirely.
# ---------------------------------------------------------------------
$rng3 = $d.Content
$found3 = $rng3.Find.Execute("This PDF version is provided under the same license.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found3) {
    $rng3.Expand($wdParagraph) | Out-Null
    $rng3.Delete()
}

# ---------------------------------------------------------------------
# 4) Remove the italic key-terms list paragraph
#    ("Uajemi, Uandishi wa kiapokaliptiki, ..."), keeping the "U"
#    Heading2 paragraph that precedes it.
# ---------------------------------------------------------------------
$rng4 = $d.Content
$found4 = $rng4.Find.Execute("Uajemi, Uandishi wa kiapokaliptiki", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found4) {
    $rng4.Expand($wdParagraph) | Out-Null
    $rng4.Delete()
}

Write-Output "done"
